$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date cell formatting (style) from A4 down to the new A5:A6 cells
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data rows
$ws.Range("A5").Value = 41551
$ws.Range("B5").Value = 3

$ws.Range("A6").Value = 41552
$ws.Range("B6").Value = 5

# Update the active selection, matching the diff
$ws.Range("B7").Select()
